$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "69.710.87"
$ws.Cells.Item(2, 5).Value = "  +2.61%  "

$ws.Cells.Item(3, 4).Value = "3.369.12"
$ws.Cells.Item(3, 5).Value = "  +3.13%  "

$ws.Cells.Item(4, 5).Value = "  +0.03%  "

$ws.Cells.Item(5, 4).Value = "'191.81"
$ws.Cells.Item(5, 5).Value = "  +3.15%  "

$ws.Cells.Item(6, 4).Value = "'591.78"
$ws.Cells.Item(6, 5).Value = "  +1.89%  "

$ws.Cells.Item(7, 5).Value = "  +1.31%  "

$ws.Cells.Item(8, 5).Value = "  +0.01%  "

$ws.Cells.Item(9, 4).Value = "'0.133"
$ws.Cells.Item(9, 5).Value = "  +2.12%  "

$ws.Cells.Item(10, 4).Value = "'6.78"
$ws.Cells.Item(10, 5).Value = "  +3.04%  "

$ws.Cells.Item(11, 5).Value = "  +1.68%  "

$ws.Cells.Item(12, 4).Value = "3.957.56"
$ws.Cells.Item(12, 5).Value = "  +3.17%  "

$ws.Cells.Item(13, 5).Value = "  -1.01%  "

$ws.Cells.Item(14, 4).Value = "'28.57"
$ws.Cells.Item(14, 5).Value = "  +3.76%  "

$ws.Cells.Item(15, 4).Value = "69.696.51"
$ws.Cells.Item(15, 5).Value = "  +2.58%  "

$ws.Cells.Item(16, 4).Value = "'0.0000171"
$ws.Cells.Item(16, 5).Value = "  +1.36%  "

$ws.Cells.Item(17, 4).Value = "3.372.89"
$ws.Cells.Item(17, 5).Value = "  +2.55%  "

$ws.Cells.Item(18, 4).Value = "'456.82"
$ws.Cells.Item(18, 5).Value = "  +14.88%  "

$ws.Cells.Item(19, 4).Value = "'5.82"
$ws.Cells.Item(19, 5).Value = "  +1.08%  "

$ws.Cells.Item(20, 4).Value = "'13.74"
$ws.Cells.Item(20, 5).Value = "  +1.42%  "

$ws.Cells.Item(21, 4).Value = "'7.89"
$ws.Cells.Item(21, 5).Value = "  +3.36%  "

$ws.Cells.Item(22, 4).Value = "'76.05"
$ws.Cells.Item(22, 5).Value = "  +6.29%  "

$ws.Cells.Item(23, 5).Value = "  -0.29%  "

$ws.Cells.Item(24, 4).Value = "3.523.26"
$ws.Cells.Item(24, 5).Value = "  +3.30%  "

$ws.Cells.Item(25, 4).Value = "'0.523"
$ws.Cells.Item(25, 5).Value = "  +2.06%  "

$ws.Cells.Item(26, 5).Value = "  +3.42%  "

$ws.Cells.Item(27, 5).Value = "  +1.86%  "

$ws.Cells.Item(28, 4).Value = "'9.46"
$ws.Cells.Item(28, 5).Value = "  -0.65%  "

$ws.Cells.Item(29, 5).Value = "  -0.33%  "

$ws.Cells.Item(30, 5).Value = "  +3.13%  "

$ws.Cells.Item(31, 4).Value = "'23.34"
$ws.Cells.Item(31, 5).Value = "  +2.84%  "

$ws.Cells.Item(32, 4).Value = "'5.56"
$ws.Cells.Item(32, 5).Value = "  +0.97%  "

$ws.Cells.Item(33, 4).Value = "'1.29"
$ws.Cells.Item(33, 5).Value = "  +2.49%  "

$ws.Cells.Item(34, 4).Value = "'6.97"
$ws.Cells.Item(34, 5).Value = "  +0.15%  "

$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 5).Value = "  +0.01%  "

$ws.Cells.Item(36, 4).Value = "'1.56"
$ws.Cells.Item(36, 5).Value = "  +6.05%  "

$ws.Cells.Item(37, 4).Value = "'164.78"
$ws.Cells.Item(37, 5).Value = "  +0.77%  "

$ws.Cells.Item(38, 5).Value = "  +2.51%  "

$ws.Cells.Item(39, 4).Value = "'27.25"
$ws.Cells.Item(39, 5).Value = "  +1.88%  "

$ws.Cells.Item(40, 4).Value = "'0.812"
$ws.Cells.Item(40, 5).Value = "  +0.36%  "

$ws.Cells.Item(41, 4).Value = "'4.60"
$ws.Cells.Item(41, 5).Value = "  +1.33%  "

$ws.Cells.Item(42, 4).Value = "'6.52"
$ws.Cells.Item(42, 5).Value = "  +2.05%  "

$ws.Cells.Item(43, 4).Value = "2.726.13"
$ws.Cells.Item(43, 5).Value = "  +2.67%  "

$ws.Cells.Item(44, 4).Value = "'2.53"
$ws.Cells.Item(44, 5).Value = "  +3.77%  "

$ws.Cells.Item(45, 4).Value = "'0.0690"
$ws.Cells.Item(45, 5).Value = "  +0.48%  "

$ws.Cells.Item(46, 4).Value = "'25.43"
$ws.Cells.Item(46, 5).Value = "  +2.35%  "

$ws.Cells.Item(47, 4).Value = "'40.97"
$ws.Cells.Item(47, 5).Value = "  +0.54%  "

$ws.Cells.Item(48, 4).Value = "'335.89"
$ws.Cells.Item(48, 5).Value = "  +0.71%  "

$ws.Cells.Item(49, 5).Value = "  +2.86%  "

$ws.Cells.Item(50, 4).Value = "'32.64"
$ws.Cells.Item(50, 5).Value = "  +5.84%  "

$ws.Cells.Item(51, 4).Value = "'1.01"
$ws.Cells.Item(51, 5).Value = "  +4.01%  "
